# Generate Report for Archive
# The handoff/handback report regenerated its rows in a new order: the
# "ab719be0-eef1-4e5f-8408-04251f37606f" item now sorts ahead of
# "ee8ee80d-703e-4c0c-bbc1-915a35bae61a" and
# "8a9aaca5-1a5a-41e5-80fa-ebad684a9799" on every sheet (Overview, zh-cn,
# de-de). Concretely this is a 3-row rotation: what used to be row 8
# becomes row 6, and the old rows 6/7 slide down to 7/8.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet (columns: A=File Name, B=zh-cn, C=de-de) ----
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A6").Value = "ab719be0-eef1-4e5f-8408-04251f37606f.md"
$ws.Range("B6").Value = "In Translation"
$ws.Range("C6").Value = "In Translation"

$ws.Range("A7").Value = "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.md"
$ws.Range("B7").Value = "In Translation"
$ws.Range("C7").Value = "In Translation"

$ws.Range("A8").Value = "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.md"
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "Ready for handoff"

# ---- zh-cn detail sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A6").Value = "ab719be0-eef1-4e5f-8408-04251f37606f.md"
$ws.Range("B6").Value = "In Translation"
$ws.Range("C6").Value = "ab719be0-eef1-4e5f-8408-04251f37606f.c643cc6f5ec3e7c167fd4a8d5ea0ea301dcfdcc4.zh-cn.xlf"
$ws.Range("D6").Value = "2016-02-17 09:29:33"

$ws.Range("A7").Value = "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.md"
$ws.Range("B7").Value = "In Translation"
$ws.Range("C7").Value = "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.05b919cc8f0ec86e6ac3a98cea0d9a70ef5e14bf.zh-cn.xlf"
$ws.Range("D7").Value = "2016-02-17 09:16:42"

$ws.Range("A8").Value = "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.md"
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.6c2925e9deb4fffac6eb59c95d2ddf7801228231.zh-cn.xlf"
$ws.Range("D8").Value = "2016-02-17 09:19:02"

# ---- de-de detail sheet ----
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A6").Value = "ab719be0-eef1-4e5f-8408-04251f37606f.md"
$ws.Range("B6").Value = "In Translation"
$ws.Range("C6").Value = "ab719be0-eef1-4e5f-8408-04251f37606f.c643cc6f5ec3e7c167fd4a8d5ea0ea301dcfdcc4.de-de.xlf"
$ws.Range("D6").Value = "2016-02-17 09:29:45"

$ws.Range("A7").Value = "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.md"
$ws.Range("B7").Value = "In Translation"
$ws.Range("C7").Value = "ee8ee80d-703e-4c0c-bbc1-915a35bae61a.05b919cc8f0ec86e6ac3a98cea0d9a70ef5e14bf.de-de.xlf"
$ws.Range("D7").Value = "2016-02-17 09:16:56"

$ws.Range("A8").Value = "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.md"
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "8a9aaca5-1a5a-41e5-80fa-ebad684a9799.6c2925e9deb4fffac6eb59c95d2ddf7801228231.de-de.xlf"
$ws.Range("D8").Value = "2016-02-17 09:19:13"
